$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared string "test" (referenced by A2) becomes "KuanWen"
$ws.Range("A2").Value = "KuanWen"

# B2 switches from the shared string "Nan" to a numeric value
$ws.Range("B2").Value = 0.06561811251427289

# C2 numeric value updated
$ws.Range("C2").Value = 5.122529692120022

# D2 keeps its shared string "Nan" (unchanged)

# E2 numeric value updated
$ws.Range("E2").Value = 0.1963057782914903
